$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.935.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.31%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.156.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.11%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.69%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.72%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.156.36"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.39%  "

# Row 11
$ws.Range("E11").Value = "  -0.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.501"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.80%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.65%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.676.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.923.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.164.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.41%  "

# Row 18
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.03%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.111"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.45%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "504.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.56%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.41%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.715"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.89%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.97%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.64%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.85%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.28%  "

# Row 28
$ws.Range("E28").Value = "  +0.13%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.30%  "

# Row 30
$ws.Range("B30").Value = "Stacks"
$ws.Range("C30").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.20%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.29%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.18%  "

# Row 33
$ws.Range("E33").Value = "  -0.11%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.43%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.66%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.54%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0893"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.19%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "480.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.46%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0416"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.93%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.75%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.992.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.95%  "

# Row 43
$ws.Range("E43").Value = "  -2.09%  "

# Row 44
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.09%  "

# Row 45
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.62%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.21%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0593"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.52%  "

# Row 49
$ws.Range("E49").Value = "  -1.38%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.61%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.99%  "
